$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1:K14").Copy($ws.Range("F1:F14"))
$ws.Columns.Item(11).Delete()
